$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '30.765.87'
$ws.Range("E2").Value = '  -0.81%  '

# Row 3
$ws.Range("D3").Value = '1.927.66'
$ws.Range("E3").Value = '  -1.41%  '

# Row 4
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.Style = "Normal"
$ws.Range("E4").Value = '  +0.04%  '

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '242.02'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.28%  '

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.03%  '

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.4857'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -0.19%  '

# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.2930'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -0.57%  '

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.06817'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -0.08%  '

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '19.12'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -0.20%  '

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '106.08'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -0.97%  '

# Row 12
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.932.27'
$ws.Range("E12").Value = '  -1.11%  '

# Row 13
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.07761'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -0.71%  '

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '5.318'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -2.50%  '

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.6957'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -0.88%  '

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '274.82'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -2.79%  '

# Row 17
$ws.Range("D17").Value = '30.753.63'
$ws.Range("E17").Value = '  -0.91%  '

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.000007658'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -0.32%  '

# Row 19
$ws.Range("E19").Value = '  +0.01%  '

# Row 20
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '12.94'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -1.84%  '

# Row 21
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '5.561'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +1.28%  '

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +0.06%  '

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '6.443'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -0.62%  '

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '9.837'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +0.30%  '

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '164.46'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -3.31%  '

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '19.43'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -2.55%  '

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '2.146'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -2.36%  '

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '0.1037'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -1.62%  '

# Row 29
$ws.Range("E29").Value = '  -2.54%  '

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.548'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -2.20%  '

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '4.345'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -2.09%  '

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.04866'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -1.27%  '

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.7558'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -0.95%  '

# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.139'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -2.63%  '

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -0.05%  '

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '2.722'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.28%  '

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.01985'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -0.98%  '

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '2.644'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -2.18%  '

# Row 40
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '6.434'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -1.35%  '

# Row 41
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '77.20'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +3.19%  '

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '2.049'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -2.48%  '

# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.8797'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -0.90%  '

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.4414'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -1.00%  '

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '107.46'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -1.63%  '

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '7.854'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -4.25%  '

# Row 47
$ws.Range("E47").Value = '  +0.03%  '

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '979.96'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -1.90%  '

# Row 49
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '36.04'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +1.04%  '

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.1234'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -1.72%  '

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '9.191'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -1.44%  '
